# Applies the crypto-price/volume refresh from the Feb 28 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '57.126.23'
$ws.Range("E2").Value = '  +1.97%  '

# Row 3
$ws.Range("D3").Value = '3.257.80'
$ws.Range("E3").Value = '  +1.26%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").Value = '''397.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.06%  '

# Row 6
$ws.Range("D6").Value = '''108.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.37%  '

# Row 7
$ws.Range("D7").Value = '''0.579'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.21%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("D9").Value = '''0.619'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.56%  '

# Row 10
$ws.Range("D10").Value = '''39.35'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.17%  '

# Row 11
$ws.Range("D11").Value = '''0.0953'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.79%  '

# Row 12
$ws.Range("E12").Value = '  +1.66%  '

# Row 13
$ws.Range("D13").Value = '3.779.42'
$ws.Range("E13").Value = '  +1.49%  '

# Row 14
$ws.Range("D14").Value = '''8.27'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.74%  '

# Row 15
$ws.Range("D15").Value = '''18.98'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.77%  '

# Row 16
$ws.Range("D16").Value = '3.260.05'
$ws.Range("E16").Value = '  +1.58%  '

# Row 17
$ws.Range("E17").Value = '  -2.04%  '

# Row 18
$ws.Range("D18").Value = '''11.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.80%  '

# Row 19
$ws.Range("D19").Value = '56.979.85'
$ws.Range("E19").Value = '  +2.11%  '

# Row 20
$ws.Range("D20").Value = '''3.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.45%  '

# Row 21
$ws.Range("D21").Value = '''0.0000108'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.94%  '

# Row 22
$ws.Range("D22").Value = '''12.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.25%  '

# Row 23
$ws.Range("D23").Value = '''293.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.89%  '

# Row 24
$ws.Range("D24").Value = '''74.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.45%  '

# Row 25
$ws.Range("E25").Value = '  -1.90%  '

# Row 26
$ws.Range("B26").Value = 'Filecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D26").Value = '''7.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.47%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''28.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.86%  '

# Row 28
$ws.Range("E28").Value = '  +0.69%  '

# Row 29
$ws.Range("D29").Value = '''7.45'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.76%  '

# Row 30
$ws.Range("E30").Value = '  -2.64%  '

# Row 31
$ws.Range("E31").Value = '  +0.06%  '

# Row 32
$ws.Range("E32").Value = '  +0.76%  '

# Row 33
$ws.Range("D33").Value = '''11.20'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.69%  '

# Row 34
$ws.Range("D34").Value = '''40.06'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +10.52%  '

# Row 35
$ws.Range("D35").Value = '''0.0489'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.52%  '

# Row 36
$ws.Range("E36").Value = '  +0.84%  '

# Row 37
$ws.Range("D37").Value = '''51.30'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.47%  '

# Row 38
$ws.Range("D38").Value = '''1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.13%  '

# Row 39
$ws.Range("E39").Value = '  -0.74%  '

# Row 40
$ws.Range("D40").Value = '''3.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.66%  '

# Row 41
$ws.Range("D41").Value = '''136.33'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.13%  '

# Row 42
$ws.Range("D42").Value = '''0.121'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.61%  '

# Row 43
$ws.Range("D43").Value = '''0.284'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.99%  '

# Row 44
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").Value = '''1.87'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.89%  '

# Row 45
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = '''3.93'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.69%  '

# Row 46
$ws.Range("D46").Value = '''16.84'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.17%  '

# Row 47
$ws.Range("D47").Value = '''22.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.24%  '

# Row 48
$ws.Range("D48").Value = '''2.20'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.16%  '

# Row 49
$ws.Range("D49").Value = '2.146.26'
$ws.Range("E49").Value = '  +0.35%  '

# Row 50
$ws.Range("D50").Value = '''2.42'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.66%  '

# Row 51
$ws.Range("D51").Value = '''1.99'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.33%  '
